# Atualizado por script em 24-11-2023 20:45
#
# This script reproduces the upstream re-scrape of the betexplorer data:
#   * Several pairs of adjacent match rows had been scraped out of order;
#     the swapped rows are restored to their correct order (everything
#     except the index/country/tournament/season/date columns, and the
#     "opening" odds timestamp columns K/O/S which are identical for both
#     matches of the pair, is exchanged between the two rows).
#   * One new match (row 160) is appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose content is swapped between the two rows of each pair.
# (A/B/C/D/E stay fixed to the row, as do K/O/S - the "opening odds"
# timestamps, which are shared between both matches of the pair.)
$swapCols = @("F","G","H","I","J","L","M","N","P","Q","R","T","U","V")

# Row-pairs that need to be swapped back into the correct order.
$rowPairs = @(
    @(7, 8),
    @(26, 27),
    @(41, 42),
    @(51, 52),
    @(99, 100),
    @(153, 154)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($col in $swapCols) {
        $addr1 = $col + $r1
        $addr2 = $col + $r2
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

# Append the newly scraped match as row 160.
$newRow = 160

# Copy number/cell formatting from the last existing data row (159) so the
# new row matches the sheet's established styling (bold/boxed index in
# column A, date format in column E).
$ws.Range("A159").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$ws.Range("E159").Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)

$ws.Range("A" + $newRow).Value = 159
$ws.Range("B" + $newRow).Value = "poland"
$ws.Range("C" + $newRow).Value = "division-2"
$ws.Range("D" + $newRow).Value = "2023-2024"
$ws.Range("E" + $newRow).Value = 45254.82291666666
$ws.Range("F" + $newRow).Value = "Olimpia Elblag"
$ws.Range("G" + $newRow).Value = 2
$ws.Range("H" + $newRow).Value = "Kotwica Kolobrzeg"
$ws.Range("I" + $newRow).Value = 3
$ws.Range("J" + $newRow).Value = 2.52
$ws.Range("K" + $newRow).Value = "23/11/2023 08:12"
$ws.Range("L" + $newRow).Value = 2.96
$ws.Range("M" + $newRow).Value = "24/11/2023 19:42"
$ws.Range("N" + $newRow).Value = 3.13
$ws.Range("O" + $newRow).Value = "23/11/2023 08:12"
$ws.Range("P" + $newRow).Value = 3.51
$ws.Range("Q" + $newRow).Value = "24/11/2023 19:42"
$ws.Range("R" + $newRow).Value = 2.52
$ws.Range("S" + $newRow).Value = "23/11/2023 08:12"
$ws.Range("T" + $newRow).Value = 2.23
$ws.Range("U" + $newRow).Value = "24/11/2023 19:42"
$ws.Range("V" + $newRow).Value = "https://www.betexplorer.com/football/poland/division-2/olimpia-elblag-kotwica-kolobrzeg/CpLaqxMD/"

Write-Host "Row swaps and append complete."
